# Apply updated cryptocurrency price/volume data to the worksheet.
# (Thu Aug  8 06:52:54 UTC 2024 GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.155.10"
$ws.Range("E2").Value = "  +0.99%  "

# Row 3
$ws.Range("D3").Value = "2.411.40"
$ws.Range("E3").Value = "  -3.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "489.22"
$ws.Range("E5").Value = "  -0.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.12"
$ws.Range("E6").Value = "  +1.38%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("E8").Value = "  +18.39%  "

# Row 9
$ws.Range("D9").Value = "2.429.28"
$ws.Range("E9").Value = "  -2.70%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.32"
$ws.Range("E10").Value = "  +9.78%  "

# Row 11
$ws.Range("E11").Value = "  +1.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  +0.56%  "

# Row 13
$ws.Range("E13").Value = "  +1.40%  "

# Row 14
$ws.Range("D14").Value = "2.839.33"
$ws.Range("E14").Value = "  -2.84%  "

# Row 15
$ws.Range("D15").Value = "57.091.98"
$ws.Range("E15").Value = "  +0.61%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.74"
$ws.Range("E16").Value = "  -2.54%  "

# Row 17
$ws.Range("E17").Value = "  -1.85%  "

# Row 18
$ws.Range("D18").Value = "2.432.87"
$ws.Range("E18").Value = "  -2.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  +4.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "323.90"
$ws.Range("E20").Value = "  +1.08%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.99"
$ws.Range("E21").Value = "  -3.07%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.98"
$ws.Range("E23").Value = "  +1.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.41"
$ws.Range("E24").Value = "  -0.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.407"
$ws.Range("E25").Value = "  -0.69%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.59%  "

# Row 27
$ws.Range("E27").Value = "  -0.45%  "

# Row 28
$ws.Range("D28").Value = "2.532.75"
$ws.Range("E28").Value = "  -2.66%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("E29").Value = "  -3.64%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0782"
$ws.Range("E30").Value = "  -2.89%  "

# Row 31
$ws.Range("E31").Value = "  +0.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.64"
$ws.Range("E32").Value = "  -0.28%  "

# Row 33
$ws.Range("E33").Value = "  +1.39%  "

# Row 34
$ws.Range("E34").Value = "  +0.56%  "

# Row 35
$ws.Range("E35").Value = "  +0.42%  "

# Row 36
$ws.Range("E36").Value = "  -0.26%  "

# Row 37
$ws.Range("E37").Value = "  -1.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("E38").Value = "  -2.14%  "

# Row 39
$ws.Range("E39").Value = "  +9.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.18"
$ws.Range("E40").Value = "  +0.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.53"
$ws.Range("E41").Value = "  +0.42%  "

# Row 42
$ws.Range("E42").Value = "  -1.31%  "

# Row 43
$ws.Range("E43").Value = "  +0.14%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.595"
$ws.Range("E44").Value = "  -3.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "268.92"
$ws.Range("E45").Value = "  +0.40%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0532"
$ws.Range("E46").Value = "  -5.63%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("E47").Value = "  +0.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0229"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49
$ws.Range("E49").Value = "  -6.29%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.878.19"
$ws.Range("E50").Value = "  -0.49%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.44"
$ws.Range("E51").Value = "  -1.61%  "
